$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 to "test" (was the long build/date string)
$ws.Range("A2").Value = "test"

# Move the active selection to A3 (it was B10)
$ws.Range("A3").Select()
